# accucor.xlsx: drop the "adductName" column from the "Corrected" sheet.
#
# The sheet's header row was:
#   Compound | C_Label | adductName | M1_mix1_T150 | M2_mix1_T150 | M3_glycerol_T150 | M4_glycerol_T150
# (columns A..G, with adductName in column C). The edit removes that
# adductName column outright (not just blanks it), which shifts the four
# sample columns left by one (D->C, E->D, F->E, G->F) and shrinks the
# sheet's used range from A1:G12 to A1:F12.
$wb = $excel.ActiveWorkbook
$corrected = $wb.Worksheets.Item("Corrected")
$corrected.Columns.Item(3).Delete()

# The workbook was also re-saved with "Corrected" as the active/selected
# tab (previously "Normalized" was active).
$corrected.Activate()
